$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> new nombre_aides (col C), new montant_total (col E)
$updates = @{
    10  = @{ C = 278209;  E = 1752512170 }
    74  = @{ C = 27992;   E = 54752531 }
    100 = @{ C = 9839;    E = 24692472 }
    117 = @{ C = 19731;   E = 56664962 }
    168 = @{ C = 285091;  E = 1212921927 }
    169 = @{ C = 562664;  E = 1286053236 }
    170 = @{ C = 367533;  E = 2847702736 }
    171 = @{ C = 115215;  E = 448609800 }
    173 = @{ C = 54396;   E = 151948108 }
    174 = @{ C = 357347;  E = 1019905253 }
    175 = @{ C = 125678;  E = 815152093 }
    177 = @{ C = 96778;   E = 174811932 }
    179 = @{ C = 235782;  E = 813286352 }
    180 = @{ C = 141523;  E = 341217178 }
    188 = @{ C = 19712;   E = 66204904 }
    255 = @{ C = 141370;  E = 414531469 }
    280 = @{ C = 95353;   E = 282552411 }
    286 = @{ C = 90609;   E = 162838526 }
    313 = @{ C = 220662;  E = 1371202181 }
    322 = @{ C = 81164;   E = 254557631 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
}
